# Update cryptos.xlsx symbol list (daily GitHub Actions refresh)
# Applies refreshed price/volume figures and re-orders two swapped coin rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.00"
$ws.Range("E2").Value = "'-0.82%"
$ws.Range("D3").Value = "'31.70"
$ws.Range("E3").Value = "'0.77%"
$ws.Range("D4").Value = "'5.087"
$ws.Range("E4").Value = "'-1.30%"
$ws.Range("D5").Value = "'0.08157"
$ws.Range("E5").Value = "'10.63%"
$ws.Range("D6").Value = "'2.580"
$ws.Range("E6").Value = "'5.96%"
$ws.Range("D7").Value = "'7.777"
$ws.Range("E7").Value = "'-1.76%"
$ws.Range("D8").Value = "'3.844"
$ws.Range("E8").Value = "'2.34%"
$ws.Range("D9").Value = "'0.9298"
$ws.Range("E9").Value = "'0.89%"
$ws.Range("D10").Value = "'0.1760"
$ws.Range("E10").Value = "'0.92%"
$ws.Range("D11").Value = "'0.07540"
$ws.Range("E11").Value = "'0.49%"
$ws.Range("D12").Value = "'0.08978"
$ws.Range("E12").Value = "'10.61%"
$ws.Range("D13").Value = "'0.03002"
$ws.Range("E13").Value = "'-1.17%"
$ws.Range("E14").Value = "'0.66%"
$ws.Range("D15").Value = "'0.001494"
$ws.Range("E15").Value = "'-0.07%"
$ws.Range("D16").Value = "'0.005751"
$ws.Range("E16").Value = "'-5.82%"
$ws.Range("D17").Value = "'3.584"
$ws.Range("E17").Value = "'3.81%"
$ws.Range("D18").Value = "'2.260"
$ws.Range("E18").Value = "'1.49%"
$ws.Range("E19").Value = "'-1.86%"
$ws.Range("D20").Value = "'0.1336"
$ws.Range("E20").Value = "'-0.11%"
$ws.Range("D21").Value = "'3.905"
$ws.Range("E21").Value = "'-16.14%"
$ws.Range("D22").Value = "'0.1697"
$ws.Range("E22").Value = "'8.12%"
$ws.Range("D23").Value = "'0.04602"
$ws.Range("E23").Value = "'-0.69%"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'1.41%"
$ws.Range("D25").Value = "'0.004468"
$ws.Range("E25").Value = "'-0.35%"
$ws.Range("E26").Value = "'-7.96%"
$ws.Range("D27").Value = "'0.0003403"
$ws.Range("E27").Value = "'81.84%"
$ws.Range("D39").Value = "'0.01772"
$ws.Range("E39").Value = "'2.69%"
$ws.Range("D40").Value = "'0.04537"
$ws.Range("E40").Value = "'0.35%"
$ws.Range("D41").Value = "'0.006929"
$ws.Range("E41").Value = "'-3.97%"
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D43").Value = "'0.002206"
$ws.Range("E43").Value = "'-0.74%"
$ws.Range("D44").Value = "'0.009712"
$ws.Range("E44").Value = "'-11.00%"
$ws.Range("D45").Value = "'0.00006420"
$ws.Range("E45").Value = "'2.02%"
$ws.Range("E46").Value = "'-0.20%"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "'0.8206"
$ws.Range("E47").Value = "'15.47%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.008739"
$ws.Range("E48").Value = "'-12.68%"
$ws.Range("D49").Value = "'0.00002096"
$ws.Range("E49").Value = "'-0.20%"
$ws.Range("D50").Value = "'0.0001996"
$ws.Range("E50").Value = "'-0.13%"

Write-Output "Updated symbol list with GitHub Actions"
